# Sensor bill-of-materials update:
#  - quantity for R2/R6 (H5) goes from 4 to 6
#  - recompute/refresh the "Total" price text (G7): $0.4043 -> $0.4067
#  - tidy up column widths for Reference/Value/Datasheet/Description
#  - leave the active selection on H9 (next blank row) as last touched in the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quantity bump for R2, R6
$ws.Range("H5").Value = 6

# Update the computed total price. The cell holds literal text (not a
# currency number), so force Text format before writing the "$..." string
# or Excel will auto-convert it to a numeric currency value; then restore
# the cell's formatting back to the sheet default.
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "$0.4067"
$ws.Range("G7").ClearFormats()

# Widen the first few columns so the BOM text isn't truncated.
$ws.Columns("A").ColumnWidth = 17.5
$ws.Columns("B").ColumnWidth = 18.333333333333332
$ws.Columns("C").ColumnWidth = 6.5
$ws.Columns("D").ColumnWidth = 45.666666666666664

# Leave the selection where the user last clicked.
[void]$ws.Range("H9").Select()
